# Apply weekly cryptos price/volume refresh (GitHub Actions scheduled update)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.963.25'
$ws.Range("E2").Value = '  +0.61%  '

$ws.Range("D3").Value = '1.640.33'
$ws.Range("E3").Value = '  +1.03%  '

$ws.Range("E4").Value = '  +0.54%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '216.00'
$ws.Range("E5").Value = '  +0.81%  '

$ws.Range("E6").Value = '  +1.53%  '

$ws.Range("E7").Value = '  +0.50%  '

$ws.Range("E8").Value = '  +0.90%  '

$ws.Range("E9").Value = '  +1.49%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.63'
$ws.Range("E10").Value = '  +0.36%  '

$ws.Range("E11").Value = '  +1.16%  '

$ws.Range("D12").Value = '1.870.34'
$ws.Range("E12").Value = '  +1.22%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.29'
$ws.Range("E13").Value = '  +1.30%  '

$ws.Range("D14").Value = '1.643.56'
$ws.Range("E14").Value = '  +1.28%  '

$ws.Range("E15").Value = '  +0.83%  '

$ws.Range("D16").Value = '0.0₃0767'
$ws.Range("E16").Value = '  +1.60%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '62.91'
$ws.Range("E17").Value = '  +0.75%  '

$ws.Range("D18").Value = '25.919.45'
$ws.Range("E18").Value = '  +0.51%  '

$ws.Range("E19").Value = '  +0.51%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '193.07'
$ws.Range("E20").Value = '  +0.57%  '

$ws.Range("E21").Value = '  +0.50%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.94'
$ws.Range("E22").Value = '  +0.21%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.25'
$ws.Range("E23").Value = '  +0.58%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.132'
$ws.Range("E24").Value = '  +6.60%  '

$ws.Range("E25").Value = '  +1.08%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '144.41'
$ws.Range("E26").Value = '  +1.93%  '

$ws.Range("E27").Value = '  +0.59%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.92'
$ws.Range("E28").Value = '  +1.25%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.53'
$ws.Range("E29").Value = '  +0.90%  '

$ws.Range("E30").Value = '  +1.03%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0499'
$ws.Range("E31").Value = '  +0.86%  '

$ws.Range("E32").Value = '  -0.91%  '

$ws.Range("E33").Value = '  +1.53%  '

$ws.Range("E34").Value = '  -2.63%  '

$ws.Range("E35").Value = '  +2.69%  '

$ws.Range("E36").Value = '  +0.47%  '

$ws.Range("D37").Value = '1.132.32'
$ws.Range("E37").Value = '  +0.69%  '

$ws.Range("E38").Value = '  -0.74%  '

$ws.Range("E39").Value = '  -0.07%  '

$ws.Range("E40").Value = '  +0.75%  '

$ws.Range("E41").Value = '  +1.84%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '99.22'

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.797'
$ws.Range("E43").Value = '  +0.45%  '

$ws.Range("D44").Value = '1.779.73'
$ws.Range("E44").Value = '  +1.20%  '

$ws.Range("E45").Value = '  +4.01%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '56.62'
$ws.Range("E46").Value = '  +0.87%  '

$ws.Range("E47").Value = '  +1.18%  '

$ws.Range("E48").Value = '  +0.44%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.76'
$ws.Range("E49").Value = '  +2.44%  '

$ws.Range("E50").Value = '  +0.17%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0959'
$ws.Range("E51").Value = '  +0.50%  '
